# Update countries & provincias Spain
# - Refresh the "last updated" timestamp in cell A1.
# - Refresh the daily COVID-19 figures for several countries (Estados Unidos,
#   Noruega, Australia, Pakistan, Serbia, and a block of African/Caribbean/
#   Central-American countries around rows 154-174). Country names / row
#   positions are unchanged - only the numeric statistics move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value2 = "Datos actualizados a 23 de Marzo de 2020 a las 22:46"

# Estados Unidos (row 6)
$ws.Range("B6").Value2 = 42751
$ws.Range("C6").Value2 = 9205
$ws.Range("E6").Value2 = 41936
$ws.Range("G6").Value2 = 101
$ws.Range("H6").Value2 = 520

# Noruega (row 17)
$ws.Range("B17").Value2 = 2621
$ws.Range("C17").Value2 = 236
$ws.Range("E17").Value2 = 2605

# Australia (row 22)
$ws.Range("B22").Value2 = 1831
$ws.Range("C22").Value2 = 222
$ws.Range("D22").Value2 = 118
$ws.Range("E22").Value2 = 1706

# Pakistan (row 32)
$ws.Range("B32").Value2 = 875
$ws.Range("C32").Value2 = 99
$ws.Range("E32").Value2 = 856

# Serbia (row 62)
$ws.Range("E62").Value2 = 243
$ws.Range("G62").Value2 = 1
$ws.Range("H62").Value2 = 3

# Gabon (row 154)
$ws.Range("C154").Value2 = 2

# Groenlandia (row 155)
$ws.Range("B155").Value2 = 5
$ws.Range("H155").Value2 = 1

# Guinea (row 156)
$ws.Range("C156").Value2 = 0

# Suazilandia (row 157)
$ws.Range("C157").Value2 = 2

# Congo (row 158)
$ws.Range("C158").Value2 = 0

# Bahamas (row 159)
$ws.Range("C159").Value2 = 1

# Namibia (row 160)
$ws.Range("C160").Value2 = 0

# Curazao (row 161)
$ws.Range("E161").Value2 = 4
$ws.Range("H161").Value2 = 0

# Republica de Yibuti (row 162)
$ws.Range("B162").Value2 = 4
$ws.Range("C162").Value2 = 1
$ws.Range("H162").Value2 = 1

# San Bartolome (row 163)
$ws.Range("C163").Value2 = 2

# Antigua y Barbuda (row 165)
$ws.Range("C165").Value2 = 0

# Cabo Verde (row 166)
$ws.Range("C166").Value2 = 2

# Zimbabue (row 173)
$ws.Range("E173").Value2 = 3
$ws.Range("G173").Value2 = 0
$ws.Range("H173").Value2 = 0

# Islas Caimanes (row 174)
$ws.Range("G174").Value2 = 1
